{"js": "// The document contains several empty \"filler\" paragraphs that were\n// justified (w:jc w:val=\"both\"). The edit removes that justification,\n// returning those paragraphs to the (unset/default) alignment while\n// leaving every other paragraph property (and the center/right aligned\n// paragraphs) untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/alignment\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.alignment === Word.Alignment.justified) {\n    paragraph.alignment = Word.Alignment.left;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains several empty \"filler\" paragraphs that were\n# justified (w:jc w:val=\"both\"). The edit removes that justification,\n# returning those paragraphs to the (unset/default) alignment while\n# leaving every other paragraph property (and the center/right aligned\n# paragraphs) untouched.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Format.Alignment -eq 3) {\n        $p.Format.Alignment = 0\n    }\n}\n"}
